# The post "「私たちが愛する言葉」" (formerly at row 606) was removed from
# the source data, so the corresponding row in the worksheet must be
# deleted. Deleting the entire row shifts every following row up by one,
# which matches the diff (old row 607 -> new row 606, ..., old row 644 ->
# new row 643) and updates the sheet's used-range dimension from
# A1:C644 to A1:C643 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(606).Delete()
